# MODIFY ENTRUST DATA IMPORT TEMPLATE [CITIC_N]
#
# The template's layout changed: a "股东代码" (shareholder code) / "交易市场"
# (trading market) / "交易类别" (trade category) set of columns was added, the
# old "合同编号" (contract no.) column was dropped, and the sample rows were
# replaced with a new day's entrusted-order export (8 rows instead of 3).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Wipe the previous sample data (old used range was A1:R4) before laying out
# the new header/data grid so no stale cells survive in columns/rows that are
# no longer populated (e.g. old column P, which held 委托类别 data that has no
# counterpart in the same position any more).
$ws.Range("A1:S8").ClearContents() | Out-Null

# ---- Header row (row 1) ----------------------------------------------
$ws.Range("A1").Value = "委托时间"
$ws.Range("B1").Value = "申请编号"
$ws.Range("C1").Value = "证券代码"
$ws.Range("D1").Value = "证券名称"
$ws.Range("E1").Value = "买卖"
$ws.Range("F1").Value = "委托类型"
$ws.Range("G1").Value = "委托状态"
$ws.Range("H1").Value = "委托价格"
$ws.Range("I1").Value = "委托数量"
$ws.Range("J1").Value = "成交价格"
$ws.Range("K1").Value = "成交数量"
$ws.Range("L1").Value = "已撤数量"
$ws.Range("M1").Value = "股东代码"
$ws.Range("N1").Value = "资金帐号"
$ws.Range("O1").Value = "交易市场"
$ws.Range("P1").Value = "返回信息"
$ws.Range("Q1").Value = "委托编号"
$ws.Range("R1").Value = "委托类别"
$ws.Range("S1").Value = "交易类别"

# ---- Data rows (rows 2-8) ---------------------------------------------
$ws.Range("A2").Value = 0.5743287037037037
$ws.Range("B2").Value = 2761
$ws.Range("C2").Value = 2798
$ws.Range("D2").Value = "帝王洁具"
$ws.Range("E2").Value = "证券买入"
$ws.Range("F2").Value = "买卖"
$ws.Range("G2").Value = "已成"
$ws.Range("H2").Value = 37.479999999999997
$ws.Range("I2").Value = 500
$ws.Range("J2").Value = 37.450000000000003
$ws.Range("K2").Value = 500
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = 211848625
$ws.Range("N2").Value = 1018000349
$ws.Range("O2").Value = "深圳"
$ws.Range("Q2").Value = 2761
$ws.Range("R2").Value = "委托"
$ws.Range("S2").Value = "波段"

$ws.Range("A3").Value = 0.57446759259259261
$ws.Range("B3").Value = 2763
$ws.Range("C3").Value = 2798
$ws.Range("D3").Value = "帝王洁具"
$ws.Range("E3").Value = "证券买入"
$ws.Range("F3").Value = "买卖"
$ws.Range("G3").Value = "已成"
$ws.Range("H3").Value = 37.479999999999997
$ws.Range("I3").Value = 800
$ws.Range("J3").Value = 37.463999999999999
$ws.Range("K3").Value = 800
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = 211848625
$ws.Range("N3").Value = 1018000349
$ws.Range("O3").Value = "深圳"
$ws.Range("Q3").Value = 2763
$ws.Range("R3").Value = "委托"
$ws.Range("S3").Value = "波段"

$ws.Range("A4").Value = 0.57550925925925933
$ws.Range("B4").Value = 2766
$ws.Range("C4").Value = 2798
$ws.Range("D4").Value = "帝王洁具"
$ws.Range("E4").Value = "证券买入"
$ws.Range("F4").Value = "买卖"
$ws.Range("G4").Value = "已成"
$ws.Range("H4").Value = 37.450000000000003
$ws.Range("I4").Value = 5700
$ws.Range("J4").Value = 37.450000000000003
$ws.Range("K4").Value = 5700
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = 211848625
$ws.Range("N4").Value = 1018000349
$ws.Range("O4").Value = "深圳"
$ws.Range("Q4").Value = 2766
$ws.Range("R4").Value = "委托"
$ws.Range("S4").Value = "波段"

$ws.Range("A5").Value = 0.57571759259259259
$ws.Range("B5").Value = 2767
$ws.Range("C5").Value = 2798
$ws.Range("D5").Value = "帝王洁具"
$ws.Range("E5").Value = "证券买入"
$ws.Range("F5").Value = "买卖"
$ws.Range("G5").Value = "已成"
$ws.Range("H5").Value = 37.6
$ws.Range("I5").Value = 20000
$ws.Range("J5").Value = 37.58
$ws.Range("K5").Value = 20000
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = 211848625
$ws.Range("N5").Value = 1018000349
$ws.Range("O5").Value = "深圳"
$ws.Range("Q5").Value = 2767
$ws.Range("R5").Value = "委托"
$ws.Range("S5").Value = "波段"

$ws.Range("A6").Value = 0.5759953703703703
$ws.Range("B6").Value = 2769
$ws.Range("C6").Value = 2798
$ws.Range("D6").Value = "帝王洁具"
$ws.Range("E6").Value = "证券买入"
$ws.Range("F6").Value = "买卖"
$ws.Range("G6").Value = "已成"
$ws.Range("H6").Value = 37.68
$ws.Range("I6").Value = 30000
$ws.Range("J6").Value = 37.68
$ws.Range("K6").Value = 30000
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = 211848625
$ws.Range("N6").Value = 1018000349
$ws.Range("O6").Value = "深圳"
$ws.Range("Q6").Value = 2769
$ws.Range("R6").Value = "委托"
$ws.Range("S6").Value = "波段"

$ws.Range("A7").Value = 0.58422453703703703
$ws.Range("B7").Value = 3008
$ws.Range("C7").Value = 2798
$ws.Range("D7").Value = "帝王洁具"
$ws.Range("E7").Value = "证券买入"
$ws.Range("F7").Value = "买卖"
$ws.Range("G7").Value = "已成"
$ws.Range("H7").Value = 37.68
$ws.Range("I7").Value = 20000
$ws.Range("J7").Value = 37.674999999999997
$ws.Range("K7").Value = 20000
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = 211848625
$ws.Range("N7").Value = 1018000349
$ws.Range("O7").Value = "深圳"
$ws.Range("Q7").Value = 3008
$ws.Range("R7").Value = "委托"
$ws.Range("S7").Value = "波段"

$ws.Range("A8").Value = 0.60089120370370364
$ws.Range("B8").Value = 3081
$ws.Range("C8").Value = 2798
$ws.Range("D8").Value = "帝王洁具"
$ws.Range("E8").Value = "证券买入"
$ws.Range("F8").Value = "买卖"
$ws.Range("G8").Value = "已成"
$ws.Range("H8").Value = 38.270000000000003
$ws.Range("I8").Value = 24300
$ws.Range("J8").Value = 38.222999999999999
$ws.Range("K8").Value = 24300
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = 211848625
$ws.Range("N8").Value = 1018000349
$ws.Range("O8").Value = "深圳"
$ws.Range("Q8").Value = 3081
$ws.Range("R8").Value = "委托"
$ws.Range("S8").Value = "波段"

# Column A holds fraction-of-day entrust times, as before - re-apply the
# same h:mm:ss time format the original rows used.
$ws.Range("A2:A8").NumberFormat = "h:mm:ss"

# Match the saved cursor position recorded in the new sheet view.
$ws.Range("H16").Select() | Out-Null
